# Updated symbol list on Tue Jan 31 06:33:54 UTC 2023 with GitHub Actions
# Refreshes Price (D) and Volume 1h % (E) columns for the affected coin rows.
# Values are prefixed with a literal apostrophe so Excel stores them as text
# (matching the workbook's existing inlineStr cells) rather than converting
# them to numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").Value = "'311.03"
$ws.Range("E2").Value = "'-1.36%"

# Row 3: OKB
$ws.Range("D3").Value = "'37.74"
$ws.Range("E3").Value = "'-3.88%"

# Row 4: HuobiToken
$ws.Range("D4").Value = "'5.074"
$ws.Range("E4").Value = "'-1.20%"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.07745"
$ws.Range("E5").Value = "'-5.18%"

# Row 6: GateToken
$ws.Range("D6").Value = "'4.353"
$ws.Range("E6").Value = "'-1.08%"

# Row 7: KuCoinToken
$ws.Range("D7").Value = "'8.199"
$ws.Range("E7").Value = "'-1.82%"

# Row 8: FTXToken
$ws.Range("D8").Value = "'1.876"
$ws.Range("E8").Value = "'-5.14%"

# Row 9: BTSEToken
$ws.Range("D9").Value = "'2.881"
$ws.Range("E9").Value = "'-11.72%"

# Row 10: MXToken
$ws.Range("D10").Value = "'0.9200"
$ws.Range("E10").Value = "'-1.66%"

# Row 11: LiechtensteinCryptoassetsExchange
$ws.Range("D11").Value = "'0.1195"
$ws.Range("E11").Value = "'-7.63%"

# Row 12: WazirX
$ws.Range("D12").Value = "'0.1911"
$ws.Range("E12").Value = "'-3.72%"

# Row 13: MandalaExchangeToken
$ws.Range("D13").Value = "'0.08896"
$ws.Range("E13").Value = "'-2.20%"

# Row 14: BitrueCoin
$ws.Range("D14").Value = "'0.03387"
$ws.Range("E14").Value = "'-4.22%"

# Row 15: BitMartToken
$ws.Range("D15").Value = "'0.09698"

# Row 16: BitForexToken
$ws.Range("D16").Value = "'0.001378"
$ws.Range("E16").Value = "'-1.79%"

# Row 17: TigerCash
$ws.Range("D17").Value = "'0.005920"
$ws.Range("E17").Value = "'-5.01%"

# Row 18: LEO
$ws.Range("E18").Value = "'-1.83%"

# Row 19: BitpandaEcosystemToken
$ws.Range("D19").Value = "'0.3406"
$ws.Range("E19").Value = "'-1.82%"

# Row 20: ProBitToken
$ws.Range("D20").Value = "'0.1280"
$ws.Range("E20").Value = "'-2.20%"

# Row 21: MCDex
$ws.Range("D21").Value = "'5.041"
$ws.Range("E21").Value = "'0.11%"

# Row 22: ZBToken
$ws.Range("D22").Value = "'0.2593"
$ws.Range("E22").Value = "'4.17%"

# Row 23: UpBots
$ws.Range("D23").Value = "'0.02106"
$ws.Range("E23").Value = "'5,595.55%"

# Row 24: CoinExToken
$ws.Range("D24").Value = "'0.04391"
$ws.Range("E24").Value = "'0.39%"

# Row 25: BitKan
$ws.Range("E25").Value = "'-2.39%"

# Row 26: HotbitToken
$ws.Range("D26").Value = "'0.004237"
$ws.Range("E26").Value = "'-10.93%"

# Row 27: NitroEx
$ws.Range("D27").Value = "'0.0001352"
$ws.Range("E27").Value = "'-65.26%"

# Row 39: One
$ws.Range("D39").Value = "'0.02094"
$ws.Range("E39").Value = "'-6.95%"

# Row 40: IDEX
$ws.Range("D40").Value = "'0.04939"
$ws.Range("E40").Value = "'-5.38%"

# Row 41: KickToken
$ws.Range("D41").Value = "'0.007862"
$ws.Range("E41").Value = "'1.49%"

# Row 42: Dexo
$ws.Range("D42").Value = "'0.009899"
$ws.Range("E42").Value = "'-3.88%"

# Row 43: BKEXToken
$ws.Range("D43").Value = "'0.1341"
$ws.Range("E43").Value = "'-4.14%"

# Row 44: CEJI
$ws.Range("D44").Value = "'0.002063"
$ws.Range("E44").Value = "'-1.80%"

# Row 45: LocalTraders
$ws.Range("D45").Value = "'0.009636"
$ws.Range("E45").Value = "'5.58%"

# Row 46: CoinLion
$ws.Range("D46").Value = "'0.00006577"
$ws.Range("E46").Value = "'-3.56%"

# Row 47: Kangarootoken
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.10%"

# Row 48: BOLO
$ws.Range("D48").Value = "'0.003046"

# Row 50: CryptobidCoin
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.10%"

# Row 51: SpecialPowerGold
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.10%"
